# 012_Registrar_Proveedor.docx
# "Agregue datos faltantes segun profe (provincia, pais, localidad, contacto)"
#
# 1) The requirements sentence about the data the system asks for when
#    registering a supplier is missing a few fields the professor asked
#    for: country, locality, province and a contact. Extend it.
# 2) Word's hidden "_GoBack" bookmark (last-edit marker) therefore moves
#    from wherever it used to sit to the spot we just edited, so re-anchor
#    it there and clear out the paragraph that used to hold it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Update the sentence with the missing data points.
# ---------------------------------------------------------------------
$oldSentence = "El sistema solicita se ingresen los siguientes datos: nombre, razón social, teléfono celular, teléfono fijo, email."
$newSentence = "El sistema solicita se ingresen los siguientes datos: nombre, razón social, teléfono celular, teléfono fijo, email, país, localidad, provincia y un contacto."

$findRange = $d.Content
$replaced = $findRange.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)
if (-not $replaced) {
    throw "Could not find the supplier-data sentence to update."
}

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark to the edit location: clear the (now
#    stale) paragraph that used to hold it, right after "El catálogo no
#    se registró.", and re-create the bookmark at the end of the
#    sentence we just edited.
# ---------------------------------------------------------------------
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $paragraph = $paragraphs.Item($i)
    if ($paragraph.Range.Text -match "cat.logo no se registr") {
        $staleBookmarkParagraph = $paragraphs.Item($i + 1)
        $staleBookmarkParagraph.Range.Delete()
        break
    }
}

$sentenceRange = $d.Content
$sentenceFound = $sentenceRange.Find.Execute("provincia y un contacto.")
if ($sentenceFound) {
    $bookmarkPos = $sentenceRange.End - 1
    $d.Range($bookmarkPos, $bookmarkPos).Bookmarks.Add("_GoBack")
}
